$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "JulianHours"

# Column widths (approximate "best fit" widths from the authored workbook)
$ws2.Columns.Item(2).ColumnWidth = 22.26
$ws2.Columns.Item(4).ColumnWidth = 12.59

# Row heights: rows 2-14 are taller (23.25), rows 15-16 use the default height
$ws2.Range("A2:F14").RowHeight = 23.25

# Header cell
$ws2.Range("D2").Value = "Gregorian Time"

# Data rows
$ws2.Range("B4").Value = 1
$ws2.Range("C4").Value = 1
$ws2.Range("D4").Value = 0.5

$ws2.Range("B5").Value = 1.5
$ws2.Range("C5").Value = 2
$ws2.Range("D5").Value = 0

# Font size 18 across the working area (A2:F14)
$ws2.Range("A2:F14").Font.Size = 18

# B column (B3:B16) originally shares the "#,##0.00000000000" number format
$ws2.Range("B3:B16").NumberFormat = "#,##0.00000000000"

# D column time-of-day values
$ws2.Range("D4:D11").NumberFormat = "[h]:mm:ss;@"

# B3:B10 switch to a shorter decimal format
$ws2.Range("B3:B10").NumberFormat = "#,##0.0000"

# Selection / active cell
$ws2.Range("D6").Select() | Out-Null
